$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020-03-30")

$ws.Range("B65").Value = -32.81588623434916
$ws.Range("C65").Value = -36.24548956597219
$ws.Range("D65").Value = -32.6366242868741
$ws.Range("E65").Value = -33.53330230258769
$ws.Range("F65").Value = -31.79406828509521
$ws.Range("G65").Value = -28.88833918996838
$ws.Range("H65").Value = -30.59865647806746
$ws.Range("I65").Value = -25.77335156301567
$ws.Range("J65").Value = -23.29390289135626
$ws.Range("K65").Value = -17.33111597263522
$ws.Range("L65").Value = -12.56852078018352
